# Generate Report for Handback
#
# - "Status" for both locales flips from "Ready for handoff" to
#   "Handed back: in sync with en-US" (Overview + per-locale sheets).
# - Each locale sheet gains a populated "Latest Target File" (F) and
#   "Latest Handback File" (G) for every data row, each a hyperlink whose
#   display text mirrors the existing Source File (A) / Latest Handoff
#   File (D) hyperlinks.
# - "Latest Handback DateTime" (H) is stamped: zh-cn -> 2016-03-24 06:37:25,
#   de-de -> 2016-03-24 06:37:33.

$wb = $excel.ActiveWorkbook

$hyperlinkColor = 15570276   # OLE (BGR) form of RGB(0x64,0x95,0xED) - matches existing hyperlink style
$underlineSingle = [Microsoft.Office.Interop.Excel.XlUnderlineStyle]::xlUnderlineStyleSingle

function Format-AsHyperlinkCell($range) {
    $range.Font.Underline = $underlineSingle
    $range.Font.Color = $hyperlinkColor
    $range.Font.Name = "Calibri"
    $range.Font.Size = 11
}

$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet: Status columns (B = zh-cn, C = de-de) ----
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("B2").Value = $newStatus
$ovw.Range("C2").Value = $newStatus
$ovw.Range("B3").Value = $newStatus
$ovw.Range("C3").Value = $newStatus

# ---- locale sheet metadata ----
# F2/F3 (Latest Target File) and G2/G3 (Latest Handback File) both mirror
# the "a.md" source / its translated xlf on every data row (matches the
# source data - row 3 gets the same handback target as row 2).
$locales = @(
    @{
        Sheet = "zh-cn"
        AUrl = "https://github.com/OpenLocalizationTest/oltest/blob/dc0aaba4907bf64c4ff0912210cfe9310a09e4b6/e2e/a.md"
        XlfDisplay = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
        XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a453b508edcfd4614b7c6e3124882799c64b86a2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
        HandbackDateTime = "2016-03-24 06:37:25"
    },
    @{
        Sheet = "de-de"
        AUrl = "https://github.com/OpenLocalizationTest/oltest/blob/dc0aaba4907bf64c4ff0912210cfe9310a09e4b6/e2e/a.md"
        XlfDisplay = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
        XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/531ec8cb65ce8aedca6f9c8b36ebec5fd14f8ad5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
        HandbackDateTime = "2016-03-24 06:37:33"
    }
)

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.Sheet)

    # Status column (C) for both data rows
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    # Row 2: Latest Target File (F2) + Latest Handback File (G2)
    $ws.Range("F2").Value = "a.md"
    Format-AsHyperlinkCell $ws.Range("F2")
    $ws.Hyperlinks.Add($ws.Range("F2"), $loc.AUrl, "", "", "a.md") | Out-Null
    Format-AsHyperlinkCell $ws.Range("F2")

    $ws.Range("G2").Value = $loc.XlfDisplay
    Format-AsHyperlinkCell $ws.Range("G2")
    $ws.Hyperlinks.Add($ws.Range("G2"), $loc.XlfUrl, "", "", $loc.XlfDisplay) | Out-Null
    Format-AsHyperlinkCell $ws.Range("G2")

    # Row 3: Latest Target File (F3) + Latest Handback File (G3)
    $ws.Range("F3").Value = "a.md"
    Format-AsHyperlinkCell $ws.Range("F3")
    $ws.Hyperlinks.Add($ws.Range("F3"), $loc.AUrl, "", "", "a.md") | Out-Null
    Format-AsHyperlinkCell $ws.Range("F3")

    $ws.Range("G3").Value = $loc.XlfDisplay
    Format-AsHyperlinkCell $ws.Range("G3")
    $ws.Hyperlinks.Add($ws.Range("G3"), $loc.XlfUrl, "", "", $loc.XlfDisplay) | Out-Null
    Format-AsHyperlinkCell $ws.Range("G3")

    # Latest Handback DateTime (H) for both rows
    $ws.Range("H2").Value = $loc.HandbackDateTime
    $ws.Range("H3").Value = $loc.HandbackDateTime
}
